$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text for B1 ("value" -> "first_release_value")
$ws.Range("B1").Value = "first_release_value"

# Rebuild the full quarter-over-quarter data series (A2:B84)
$data = New-Object 'object[,]' 83,2
$data[0,0] = 38398
$data[0,1] = 0.4001302730732021
$data[1,0] = 38487
$data[1,1] = 1.534309226294653
$data[2,0] = 38579
$data[2,1] = -0.2158762754026498
$data[3,0] = 38671
$data[3,1] = -0.1573369388209471
$data[4,0] = 38763
$data[4,1] = 0.5441785023706558
$data[5,0] = 38852
$data[5,1] = -0.1273572285275435
$data[6,0] = 38944
$data[6,1] = 0.510068525034896
$data[7,0] = 39036
$data[7,1] = -0.009760955203091726
$data[8,0] = 39128
$data[8,1] = 1.786050651751793
$data[9,0] = 39217
$data[9,1] = -0.03834288659695062
$data[10,0] = 39309
$data[10,1] = -0.0575559711994913
$data[11,0] = 39401
$data[11,1] = -0.6200525020039009
$data[12,0] = 39493
$data[12,1] = 1.257433230729447
$data[13,0] = 39583
$data[13,1] = 0.8837904892317567
$data[14,0] = 39675
$data[14,1] = 0.8108064919307481
$data[15,0] = 39767
$data[15,1] = -0.07476001263452758
$data[16,0] = 39859
$data[16,1] = 0.2524978494830066
$data[17,0] = 39948
$data[17,1] = 0.5778148852415939
$data[18,0] = 40040
$data[18,1] = 0.1185568564730346
$data[19,0] = 40132
$data[19,1] = -0.1977656654399595
$data[20,0] = 40224
$data[20,1] = 1.071871022829441
$data[21,0] = 40313
$data[21,1] = -1.091011900795806
$data[22,0] = 40405
$data[22,1] = 1.120967691003898
$data[23,0] = 40497
$data[23,1] = -0.1775928823643795
$data[24,0] = 40589
$data[24,1] = 1.28981182300268
$data[25,0] = 40678
$data[25,1] = 0.5923450763659872
$data[26,0] = 40770
$data[26,1] = 0.5531759638372762
$data[27,0] = 40862
$data[27,1] = 0.5788603179058356
$data[28,0] = 40954
$data[28,1] = 0.1947850960503388
$data[29,0] = 41044
$data[29,1] = -0.2481858862331165
$data[30,0] = 41136
$data[30,1] = 0.3732050716642448
$data[31,0] = 41228
$data[31,1] = 0.141299961337424
$data[32,0] = 41320
$data[32,1] = -0.1411005862636046
$data[33,0] = 41409
$data[33,1] = -0.2128461555332564
$data[34,0] = 41501
$data[34,1] = 0.4532479246724535
$data[35,0] = 41593
$data[35,1] = -0.3446087745608111
$data[36,0] = 41685
$data[36,1] = 0.4255979180752121
$data[37,0] = 41774
$data[37,1] = 0.4461687925667093
$data[38,0] = 41866
$data[38,1] = 0.6373066379050414
$data[39,0] = 41958
$data[39,1] = 0.3261422475203943
$data[40,0] = 42050
$data[40,1] = 0.6788370390783598
$data[41,0] = 42139
$data[41,1] = 0.6601374471387373
$data[42,0] = 42231
$data[42,1] = 1.255382587579845
$data[43,0] = 42323
$data[43,1] = 0.8728685839363095
$data[44,0] = 42415
$data[44,1] = 0.4694885089849095
$data[45,0] = 42505
$data[45,1] = 1.152137745180852
$data[46,0] = 42597
$data[46,1] = 0.9596379771730028
$data[47,0] = 42689
$data[47,1] = 0.2682953781150843
$data[48,0] = 42781
$data[48,1] = 0.4191917022489378
$data[49,0] = 42870
$data[49,1] = 0.1682050168937224
$data[50,0] = 42962
$data[50,1] = -0.03534872415686152
$data[51,0] = 43054
$data[51,1] = 0.5002605909365485
$data[52,0] = 43146
$data[52,1] = -0.5240674734835977
$data[53,0] = 43235
$data[53,1] = 0.5854015665873362
$data[54,0] = 43327
$data[54,1] = 0.2
$data[55,0] = 43419
$data[55,1] = 1.628071843823122
$data[56,0] = 43511
$data[56,1] = -0.3
$data[57,0] = 43600
$data[57,1] = 0.5022917647287812
$data[58,0] = 43692
$data[58,1] = 0.754883892913071
$data[59,0] = 43784
$data[59,1] = 0.3494637214130449
$data[60,0] = 43876
$data[60,1] = 0.2
$data[61,0] = 43966
$data[61,1] = 1.450185044412038
$data[62,0] = 44058
$data[62,1] = 0.3499999990000049
$data[63,0] = 44150
$data[63,1] = 0.7000000000000171
$data[64,0] = 44242
$data[64,1] = 0.7000000000000171
$data[65,0] = 44331
$data[65,1] = 0.7999999999999972
$data[66,0] = 44423
$data[66,1] = -0.7999999999999972
$data[67,0] = 44515
$data[67,1] = 0.4999999999999858
$data[68,0] = 44607
$data[68,1] = 0.00000000000002842170943040401
$data[69,0] = 44696
$data[69,1] = 0.5000000000000142
$data[70,0] = 44788
$data[70,1] = 0.7000000000000171
$data[71,0] = 44880
$data[71,1] = -0.9999999999999858
$data[72,0] = 44972
$data[72,1] = -0.5
$data[73,0] = 45061
$data[73,1] = 2.799999999999997
$data[74,0] = 45153
$data[74,1] = 0
$data[75,0] = 45245
$data[75,1] = 0.2000000000000028
$data[76,0] = 45337
$data[76,1] = -0.09999999999999432
$data[77,0] = 45427
$data[77,1] = 0.2000000000000028
$data[78,0] = 45519
$data[78,1] = 0.09999999999999432
$data[79,0] = 45611
$data[79,1] = 0.09999999999999432
$data[80,0] = 45703
$data[80,1] = 0
$data[81,0] = 45792
$data[81,1] = 0.4999999999999858
$data[82,0] = 45884
$data[82,1] = 0.4000000000000199

$ws.Range("A2:B84").Value = $data

# Newly added rows (beyond the original A1:B53 range) need the same
# date-formatted / bordered / bold-centered style as the existing A column cells.
$newA = $ws.Range("A54:A84")
$newA.NumberFormat = "YYYY-MM-DD HH:MM:SS"
$newA.HorizontalAlignment = -4108
$newA.VerticalAlignment = -4160
$newA.Font.Bold = $true
$newA.Borders.LineStyle = 1

